$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "KEY ACHIEVEMENTS AND IMPACT" section and the first Heading3
# ("Software Development and Innovation") beneath it, so the script is
# resilient to this section's absolute paragraph index.
# ---------------------------------------------------------------------------
$startIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Software Development and Innovation") {
        $startIdx = $i
        break
    }
}
if ($startIdx -eq -1) {
    throw "Could not locate 'Software Development and Innovation' heading paragraph"
}

# Step A: retitle the first Heading3 of this section.
$pHeading1 = $d.Paragraphs.Item($startIdx)
$pHeading1.Range.Text = "Technical Innovation & Platform Development"

# Step B: strengthen the first bullet under it.
$idx = $startIdx + 1
$pBullet1 = $d.Paragraphs.Item($idx)
$pBullet1.Range.Text = "• Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide"

# Step C: insert 4 new bullets right after it (before the unchanged
# "Developed boundary estimation..." bullet).
$newTechBullets = @(
    "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party",
    "• Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors",
    "• Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls",
    "• Created SimCrisis platform for humanitarian intervention modeling used by International Red Cross and UNICEF"
)
foreach ($bulletText in $newTechBullets) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newP = $d.Paragraphs.Item($idx)
    $newP.Range.Text = $bulletText
}

# $idx now points at the last inserted bullet ("Created SimCrisis...").
# The next paragraph (unchanged) is "• Developed boundary estimation...".
$idx = $idx + 1  # "Developed boundary estimation..." (left untouched)

# The paragraph after that is "• Created econometric simulation platform..."
# -- this is the one that gets converted into the "Data Engineering &
# Analytics" Heading3, so capture its index before mutating anything.
$idx = $idx + 1
$econIdx = $idx

# Step D: insert the 6 "Data Engineering & Analytics" bullets right after
# the still-Normal-styled "Created econometric simulation..." paragraph, so
# the new bullets naturally inherit Normal style (no explicit style needed).
$newDataBullets = @(
    "• Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change",
    "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%",
    "• Developed advanced data pipelines for machine learning applications enhancing consumer segmentation and predictive modeling",
    "• Built fraud detection systems for campaign finance data analysis across multi-terabyte datasets",
    "• Transformed small data team into big data engineering team using Hadoop Clusters and Hive on AWS",
    "• Introduced version control and Agile methodologies, improving project delivery timelines by 40%"
)
$cursor = $econIdx
foreach ($bulletText in $newDataBullets) {
    $p = $d.Paragraphs.Item($cursor)
    $p.Range.InsertParagraphAfter()
    $cursor = $cursor + 1
    $newP = $d.Paragraphs.Item($cursor)
    $newP.Range.Text = $bulletText
}

# Step E: now convert the original "Created econometric simulation..."
# paragraph into the new Heading3 "Data Engineering & Analytics".
$pEcon = $d.Paragraphs.Item($econIdx)
$pEcon.Style = "Heading 3"
$pEcon.Range.Text = "Data Engineering & Analytics"

# $cursor currently points at the last inserted Data-Engineering bullet
# ("Introduced version control and Agile methodologies...").
# Step F: insert the new "Research Leadership & Client Success" Heading3
# right after it.
$p = $d.Paragraphs.Item($cursor)
$p.Range.InsertParagraphAfter()
$cursor = $cursor + 1
$pHeading3 = $d.Paragraphs.Item($cursor)
$pHeading3.Style = "Heading 3"
$pHeading3.Range.Text = "Research Leadership & Client Success"

# The paragraph right after the new heading is the unchanged
# "• Built comprehensive survey operations platform from RFP through
# deployment" bullet. Locate it explicitly (it is immediately next).
$surveyIdx = $cursor + 1

# Step G: insert the 3 "Research Leadership & Client Success" bullets
# BEFORE the unchanged survey bullet, so they inherit Normal style cleanly.
$newResearchBullets = @(
    "• Led multi-million dollar research projects involving sensitive consumer data with privacy compliance",
    "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders",
    "• Delivered actionable consumer insights and market intelligence for political candidates and major organizations"
)
foreach ($bulletText in $newResearchBullets) {
    $p = $d.Paragraphs.Item($surveyIdx)
    $p.Range.InsertParagraphBefore()
    $newP = $d.Paragraphs.Item($surveyIdx)
    $newP.Range.Text = $bulletText
    $surveyIdx = $surveyIdx + 1
}

# $surveyIdx now points at the unchanged survey bullet itself.
# Step H: append 2 more bullets right after it.
$newClosingBullets = @(
    "• Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership",
    "• Redistricting analysis used in court cases with rigorous methodology and expert testimony"
)
$cursor = $surveyIdx
foreach ($bulletText in $newClosingBullets) {
    $p = $d.Paragraphs.Item($cursor)
    $p.Range.InsertParagraphAfter()
    $cursor = $cursor + 1
    $newP = $d.Paragraphs.Item($cursor)
    $newP.Range.Text = $bulletText
}
